$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Strip the "id" prefix from the numeric identifiers in the feature labels ---
$ws.Range("C4").Value = "Численность - popsize (чел.) (8112027)"
$ws.Range("D4").Value = "Ср. числ. работн. орг. -  avgemployers (чел.) (8123005)"
$ws.Range("E4").Value = "Площ. торг. зал. маг. - shoparea (кв.м.) (8002002)"

$ws.Range("C5").Value = "Миг. сальдо - saldo (чел.) (8112021 - 8112022)"
$ws.Range("D5").Value = "Сред. зп. - avgsalary (руб.) (8123007)"
$ws.Range("E5").Value = "Кол-во мест в рест,каф,бар - foodseats (место) (8002004)"

$ws.Range("E6").Value = "Обор. роз. (кроме авто.) - retailturnover (тыс. руб.) (8201003)"
$ws.Range("E7").Value = "Обор. Общепит - foodservturnover (тыс. руб.) (8201006)"

# --- Second block: header text + ids without "id" prefix ---
$ws.Range("D10").Value = "Уровень жизни (разное + новое)"

$ws.Range("C11").Value = "Введ. жил. дом. - consnewareas (кв. м.) (8010001)"
$ws.Range("D11").Value = "Жил. площ.на одного чел. - livarea (кв. м) (8211001)"

# --- New feature: livestock, added for "Сельское хозяйство" column ---
$ws.Range("E11").Value = "Поголовье скота - livestock (сум. всех видов, шт.) (8007020)"

$ws.Range("C12").Value = "Введ. кварт. - consnewapt (шт. на 1000 чел.) (8215002)"
$ws.Range("D12").Value = "Число спорт. сооруж. - sportsvenue (шт.) (8003001)"

$ws.Range("D13").Value = "Объекты быт. обслу. - servicesnum (шт.) (8001001 & 8401011)"

$ws.Range("D14").Value = "Длина дорог - roadslen (км) (8006005)"

# --- Column E width shrinks slightly ---
$ws.Columns("E").ColumnWidth = 61.5

# --- Active selection moves to H11 ---
$ws.Range("H11").Select()
